$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "mode"
$ws.Range("B10").Value = "p"

$ws.Range("B11").Select()
